$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force them to remain plain text (matching the source workbook, where every
# data cell -- numeric-looking or not -- is stored as a literal string).
$ws.Range("D2").Value = '63.258.76'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.553.13'
$ws.Range("E3").Value = '  +5.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.61'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.24'
$ws.Range("E6").Value = '  +5.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("D9").Value = '2.551.52'
$ws.Range("E9").Value = '  +5.06%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.34'
$ws.Range("E14").Value = '  +4.29%  '
$ws.Range("D15").Value = '3.011.69'
$ws.Range("E15").Value = '  +5.07%  '
$ws.Range("D16").Value = '63.168.86'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = '2.553.93'
$ws.Range("E18").Value = '  +5.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.53'
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '336.23'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.29'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  -2.55%  '
$ws.Range("E26").Value = '  +4.46%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.42'
$ws.Range("E28").Value = '  +3.23%  '
$ws.Range("B29").Value = 'SuiNetwork'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.47'
$ws.Range("E29").Value = '  +10.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.17'
$ws.Range("E30").Value = '  +8.26%  '
$ws.Range("D31").Value = '0.0₃0815'
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.44'
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  +7.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '415.18'
$ws.Range("E35").Value = '  +12.22%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.88'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.40'
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("E40").Value = '  +3.88%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.41'
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '152.28'
$ws.Range("E43").Value = '  +2.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.74'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.69'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.611'
$ws.Range("E46").Value = '  +3.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0967'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0522'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  +5.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.42'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.78'
$ws.Range("E51").Value = '  +3.26%  '
